# Apply updated NATMI Vegfb-Flt1 values for columns E, G-K, M-T across rows 2-37.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Row number, then the new values for columns E, G, H, I, J, K, M, N, O, P, Q, R, S, T (in that order).
$rowData = @(
  ,@(2, 2, 1.3034235, 2.606847, 0.05019481880249994, 0.03724108001642075, 2, 117.8215865, 235.643173, 0.532002292588051, 0.4324630118342522, 153.5714246513828, 614.285698605531, 0.02670375867897178, 0.0161053896278617)
  ,@(3, 2, 1.3034235, 2.606847, 0.05019481880249994, 0.03724108001642075, 3, 0.235393, 0.706179, 0.00106287497375686, 0.001296011649079689, 0.3068167679355, 1.840900607613, 0.00005335081671743747, 0.00004826487352559011)
  ,@(4, 2, 1.3034235, 2.606847, 0.05019481880249994, 0.03724108001642075, 3, 0.5852593333333332, 1.755778, 0.002642633802014605, 0.003222283218840886, 0.7628407686609999, 4.577044611966, 0.0001326465248534846, 0.0001200013071884232)
  ,@(5, 2, 1.3034235, 2.606847, 0.05019481880249994, 0.03724108001642075, 3, 2.51651, 7.549530000000001, 0.01136285063790714, 0.01385523900466678, 3.280078271985, 19.68046963191, 0.0005703562288496199, 0.0005159840644194295)
  ,@(6, 2, 1.3034235, 2.606847, 0.05019481880249994, 0.03724108001642075, 3, 98.61275733333333, 295.838272, 0.4452682614298567, 0.5429357808085564, 128.534185308064, 771.205111848384, 0.02235015970097583, 0.02021951485686933)
  ,@(7, 2, 1.3034235, 2.606847, 0.05019481880249994, 0.03724108001642075, 2, 1.696687, 3.393374, 0.007661086568413697, 0.006227673484603961, 2.2115017079445, 8.846006831778, 0.0003845468521317916, 0.0002319252865562779)
  ,@(8, 3, 4.086202333333333, 12.258607, 0.1573595885849988, 0.1751248785129528, 2, 117.8215865, 235.643173, 0.532002292588051, 0.4324630118342522, 481.4428416733352, 2888.657050040011, 0.08371566188793184, 0.0757350324088191)
  ,@(9, 3, 4.086202333333333, 12.258607, 0.1573595885849988, 0.1751248785129528, 3, 0.235393, 0.706179, 0.00106287497375686, 0.001296011649079689, 0.9618634258503334, 8.656770832653, 0.0001672535685876708, 0.0002269638825964522)
  ,@(10, 3, 4.086202333333333, 12.258607, 0.1573595885849988, 0.1751248785129528, 3, 0.5852593333333332, 1.755778, 0.002642633802014605, 0.003222283218840886, 2.391488053471778, 21.523392481246, 0.0004158437678658293, 0.0005643019572338367)
  ,@(11, 3, 4.086202333333333, 12.258607, 0.1573595885849988, 0.1751248785129528, 3, 2.51651, 7.549530000000001, 0.01136285063790714, 0.01385523900466678, 10.28296903385667, 92.54672130471002, 0.001788053501533858, 0.002426397047460196)
  ,@(12, 3, 4.086202333333333, 12.258607, 0.1573595885849988, 0.1751248785129528, 3, 98.61275733333333, 295.838272, 0.4452682614298567, 0.5429357808085564, 402.9516791119005, 3626.565112007105, 0.07006723042855993, 0.09508156265443363)
  ,@(13, 3, 4.086202333333333, 12.258607, 0.1573595885849988, 0.1751248785129528, 2, 1.696687, 3.393374, 0.007661086568413697, 0.006227673484603961, 6.933006378336334, 41.59803827001801, 0.001205545430519639, 0.001090620562409606)
  ,@(14, 3, 5.424432333333333, 16.273297, 0.2088948051635471, 0.232478222046779, 2, 117.8215865, 235.643173, 0.532002292588051, 0.4324630118342522, 639.1152233752301, 3834.691340251381, 0.1111325152567413, 0.1005382320922221)
  ,@(15, 3, 5.424432333333333, 16.273297, 0.2088948051635471, 0.232478222046779, 3, 0.235393, 0.706179, 0.00106287497375686, 0.001296011649079689, 1.276873400240333, 11.491860602163, 0.0002220290605561495, 0.0003012944839299602)
  ,@(16, 3, 5.424432333333333, 16.273297, 0.2088948051635471, 0.232478222046779, 3, 0.5852593333333332, 1.755778, 0.002642633802014605, 0.003222283218840886, 3.174699651118444, 28.572296860066, 0.0005520324731904446, 0.0007491106736473012)
  ,@(17, 3, 5.424432333333333, 16.273297, 0.2088948051635471, 0.232478222046779, 3, 2.51651, 7.549530000000001, 0.01136285063790714, 0.01385523900466678, 13.65063821115667, 122.85574390041, 0.002373640470108099, 0.003221041329838118)
  ,@(18, 3, 5.424432333333333, 16.273297, 0.2088948051635471, 0.232478222046779, 3, 98.61275733333333, 295.838272, 0.4452682614298567, 0.5429357808085564, 534.918229358087, 4814.264064222784, 0.09301422671690127, 0.1262207450079529)
  ,@(19, 3, 5.424432333333333, 16.273297, 0.2088948051635471, 0.232478222046779, 2, 1.696687, 3.393374, 0.007661086568413697, 0.006227673484603961, 9.203563822346332, 55.221382934078, 0.001600361186049847, 0.001447798459188597)
  ,@(20, 3, 6.646854666666667, 19.940564, 0.2559702702919538, 0.2848683254124845, 2, 117.8215865, 235.643173, 0.532002292588051, 0.4324630118342522, 783.1429620615953, 4698.857772369573, 0.1361767706297025, 0.1231950139840629)
  ,@(21, 3, 6.646854666666667, 19.940564, 0.2559702702919538, 0.2848683254124845, 3, 0.235393, 0.706179, 0.00106287497375686, 0.001296011649079689, 1.564623060550667, 14.081607544956, 0.0002720643943190968, 0.0003691926681884035)
  ,@(22, 3, 6.646854666666667, 19.940564, 0.2559702702919538, 0.2848683254124845, 3, 0.5852593333333332, 1.755778, 0.002642633802014605, 0.003222283218840886, 3.890133730976888, 35.011203578792, 0.0006764356885843321, 0.0009179264245559536)
  ,@(23, 3, 6.646854666666667, 19.940564, 0.2559702702919538, 0.2848683254124845, 3, 2.51651, 7.549530000000001, 0.01136285063790714, 0.01385523900466678, 16.72687623721334, 150.54188613492, 0.002908551949072191, 0.003946918733449166)
  ,@(24, 3, 6.646854666666667, 19.940564, 0.2559702702919538, 0.2848683254124845, 3, 98.61275733333333, 295.838272, 0.4452682614298567, 0.5429357808085564, 655.4646662739343, 5899.181996465409, 0.1139754372306288, 0.1546652066854532)
  ,@(25, 3, 6.646854666666667, 19.940564, 0.2559702702919538, 0.2848683254124845, 2, 1.696687, 3.393374, 0.007661086568413697, 0.006227673484603961, 11.27763190382267, 67.66579142293601, 0.001961010399646911, 0.001774066916774863)
  ,@(26, 3, 1.907159333333333, 5.721477999999999, 0.07344467639578636, 0.08173629676394162, 2, 117.8215865, 235.643173, 0.532002292588051, 0.4324630118342522, 224.7045383616156, 1348.227230169694, 0.03907273622094586, 0.03534792507471243)
  ,@(27, 3, 1.907159333333333, 5.721477999999999, 0.07344467639578636, 0.08173629676394162, 3, 0.235393, 0.706179, 0.00106287497375686, 0.001296011649079689, 0.4489319569513333, 4.040387612561999, 0.0000780625084967525, 0.0001059311927587028)
  ,@(28, 3, 1.907159333333333, 5.721477999999999, 0.07344467639578636, 0.08173629676394162, 3, 0.5852593333333332, 1.755778, 0.002642633802014605, 0.003222283218840886, 1.116182799987111, 10.045645199884, 0.0001940873844215292, 0.0002633774974326477)
  ,@(29, 3, 1.907159333333333, 5.721477999999999, 0.07344467639578636, 0.08173629676394162, 3, 2.51651, 7.549530000000001, 0.01136285063790714, 0.01385523900466678, 4.799385533926666, 43.19446980534, 0.0008345408880347446, 0.001132475927020783)
  ,@(30, 3, 1.907159333333333, 5.721477999999999, 0.07344467639578636, 0.08173629676394162, 3, 98.61275733333333, 295.838272, 0.4452682614298567, 0.5429357808085564, 188.0702405340018, 1692.632164806016, 0.03270258337003023, 0.04437756010393053)
  ,@(31, 3, 1.907159333333333, 5.721477999999999, 0.07344467639578636, 0.08173629676394162, 2, 1.696687, 3.393374, 0.007661086568413697, 0.006227673484603961, 3.235852447795333, 19.415114686772, 0.0005626660238572494, 0.0005090269680865197)
  ,@(32, 2, 6.5992195, 13.198439, 0.254135840761214, 0.1885511972474212, 2, 117.8215865, 235.643173, 0.532002292588051, 0.4324630118342522, 777.5305111517367, 3110.122044606947, 0.1352008499137577, 0.08154141864657395)
  ,@(33, 2, 6.5992195, 13.198439, 0.254135840761214, 0.1885511972474212, 3, 0.235393, 0.706179, 0.00106287497375686, 0.001296011649079689, 1.5534100757635, 9.320460454581001, 0.0002701146250797529, 0.0002443645480805801)
  ,@(34, 2, 6.5992195, 13.198439, 0.254135840761214, 0.1885511972474212, 3, 0.5852593333333332, 1.755778, 0.002642633802014605, 0.003222283218840886, 3.862254805090333, 23.173528830542, 0.0006715879630989853, 0.0006075653587827233)
  ,@(35, 2, 6.5992195, 13.198439, 0.254135840761214, 0.1885511972474212, 3, 2.51651, 7.549530000000001, 0.01136285063790714, 0.01385523900466678, 16.607001863945, 99.64201118367001, 0.002887707600308629, 0.002612421902479091)
  ,@(36, 2, 6.5992195, 13.198439, 0.254135840761214, 0.1885511972474212, 3, 98.61275733333333, 295.838272, 0.4452682614298567, 0.5429357808085564, 650.7672311429013, 3904.603386857408, 0.1131586239827607, 0.1023711914999168)
  ,@(37, 2, 6.5992195, 13.198439, 0.254135840761214, 0.1885511972474212, 2, 1.696687, 3.393374, 0.007661086568413697, 0.006227673484603961, 11.1968099357965, 44.78723974318601, 0.001946956676208259, 0.001174235291588096)
)

$targetCols = @(5, 7, 8, 9, 10, 11, 13, 14, 15, 16, 17, 18, 19, 20)

foreach ($entry in $rowData) {
  $rowNum = $entry[0]
  for ($i = 0; $i -lt $targetCols.Count; $i++) {
    $colNum = $targetCols[$i]
    $value = $entry[$i + 1]
    $ws.Cells.Item($rowNum, $colNum).Value = $value
  }
}
